$wb = $excel.ActiveWorkbook

# Work on the "Repayment Schedule" sheet: insert a new blank column before
# column N. This pushes the old "Late" and "Outstanding" columns one place
# to the right and creates a blank spacer column.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N:N").Insert()
$wsSchedule.Columns("N:N").ColumnWidth = 10.3

# Update the selection on the Repayment Schedule sheet.
$wsSchedule.Range("R8").Select()

# Update the selection on the Transactions sheet (kept at B2, but it is no
# longer the active/selected tab).
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("B2").Select()

# Make "Repayment Schedule" the active sheet/tab (it was "Transactions").
$wsSchedule.Activate()
$wsSchedule.Range("R8").Select()
